$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Golang Architect", "https://www.dice.com/job-detail/a5af2b1e-5691-4ff6-86fc-51996fdb4cd5", "Remote", "Third Party, Contract", "Depends on Experience", "InfiCare Technologies"),
    @("Golang Architect / Principal Backend Architect", "https://www.dice.com/job-detail/8d4bdef4-f002-49a4-9040-f41eb8cba37a", "Atlanta, Georgia", "Third Party", "Depends on Experience", "Keylent"),
    @("Go Lang Developer", "https://www.dice.com/job-detail/093799d8-8c91-4d22-a714-54460b7d10e0", "Remote", "Contract, Third Party", "Depends on Experience", "Elista Global LLC")
)

$startRow = 49
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
}
